# Apply updated dSF (column F) values for specific rows on Sheet1.
# These rows had their dSF figure "repulled" from the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -10
    3  = -2
    5  = -5
    6  = 1
    9  = -8
    13 = 6
    16 = -1
    23 = -2
    31 = 0
    32 = -3
    34 = -8
    37 = 4
    38 = -4
    39 = 4
    44 = -7
    47 = -12
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

$wb.Save()
